$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Repeating 7-value cycle used to populate column A
$cycle = @(3174466432, 3247439861, 3104023154, 3215996243, 3164602900, 3162924200, 3125278534)

$totalRows = 336
$values = New-Object 'object[,]' $totalRows, 1
for ($i = 0; $i -lt $totalRows; $i++) {
    $values[$i, 0] = $cycle[$i % 7]
}

$ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($totalRows, 1)).Value = $values

# Update selection / view to match new data extent
$ws.Range("A337").Select()
$excel.ActiveWindow.ScrollRow = 313
